$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145:224 down to 146:225
$ws.Rows("145:145").Insert(-4121)  # xlShiftDown = -4121

# Populate the newly inserted row 145 with the new data record
$ws.Range("A145").Value2 = 4
$ws.Range("B145").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C145").Value2 = "Los Lagos"
$ws.Range("D145").Value2 = 44452
$ws.Range("E145").Value2 = 10
$ws.Range("F145").Value2 = 100114001
$ws.Range("G145").Value2 = "Papa"
$ws.Range("H145").Value2 = "Asterix"
$ws.Range("I145").Value2 = "1a (guarda)"
$ws.Range("J145").Value2 = 250
$ws.Range("K145").Value2 = 7000
$ws.Range("L145").Value2 = 7500
$ws.Range("M145").Value2 = 7200
$ws.Range("N145").Value2 = "$/saco 25 kilos"
$ws.Range("O145").Value2 = "Provincia de Llanquihue"
$ws.Range("P145").Value2 = 288
$ws.Range("Q145").Value2 = 25
$ws.Range("R145").Value2 = "Hortaliza"
